# Actualizacion de la base de datos del Estado de Cuenta:
# se eliminan los trabajadores anteriores y se agregan los nuevos
# (manteniendo el formato/estilo existente de cada fila).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fila 16: ESMERALDA MARIA RAMOS LOPEZ -> RONALD PALENCIA CABARCAS
$ws.Range("C16").Value = "9147189"
$ws.Range("D16").Value = "RONALD PALENCIA CABARCAS"
$ws.Range("G16").Value = 738000

# Fila 17: FRANCISCO ANTONIO LLERENA AGUILAR -> DAVID ANTONIO ANGULO ORTIZ
$ws.Range("C17").Value = "1047385679"
$ws.Range("D17").Value = "DAVID ANTONIO ANGULO ORTIZ"
$ws.Range("G17").Value = 738000

# Fila 18: RONALD PALENCIA CABARCAS -> FRANCISCO ANTONIO LLERENA AGUILAR
$ws.Range("C18").Value = "1050949909"
$ws.Range("D18").Value = "FRANCISCO ANTONIO LLERENA AGUILAR"
$ws.Range("G18").Value = 781242

# Fila 19: SUSANA PATRICIA NOEL PEREZ -> HANNY MARGARITA VELASQUEZ ARELLANO
$ws.Range("C19").Value = "1050945732"
$ws.Range("D19").Value = "HANNY MARGARITA VELASQUEZ ARELLANO"
$ws.Range("G19").Value = 738000

# Fila 20: DAVID ANTONIO ANGULO ORTIZ -> ESMERALDA MARIA RAMOS LOPEZ
$ws.Range("C20").Value = "26007512"
$ws.Range("D20").Value = "ESMERALDA MARIA RAMOS LOPEZ"
$ws.Range("G20").Value = 738000

# Fila 21: HANNY MARGARITA VELASQUEZ ARELLANO -> SUSANA PATRICIA NOEL PEREZ
$ws.Range("C21").Value = "32907043"
$ws.Range("D21").Value = "SUSANA PATRICIA NOEL PEREZ"
$ws.Range("G21").Value = 738000
